$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("anagrafica_aziendale")

# --- 1. Fill in the four "incentivo" formula columns (H:K) for rows 2-29 ---
for ($r = 2; $r -le 29; $r++) {
    $ws.Range("H$r").Formula = "=IF(G$r>=10,100,0)"
    $ws.Range("I$r").Formula = "=IF(D$r=""Produzione"",100,0)"
    $ws.Range("J$r").Formula = "=IF(AND(D$r=""Amministrazione"",G$r>=10),100,0)"
    $ws.Range("K$r").Formula = "=IF(OR(D$r=""Direzione"",D$r=""Commerciale""),100,0)"
}

# --- 2. Match the formatting of the first data row (H2:K2) onto the rest (H3:K29) ---
$null = $ws.Range("H2:K2").Copy()
$null = $ws.Range("H3:K29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Move / resize the text box shape ---
$shp = $ws.Shapes.Item("CasellaDiTesto 1")
$shp.Left = 1013.0
$shp.Top = 16.0
$shp.Width = 220.0
$shp.Height = 75.0

# --- 4. Update the active selection on the sheet ---
$null = $ws.Range("I3").Select()
